$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 37, shifting existing rows 37-38 down to 38-39.
$ws.Rows(37).Insert()

# Populate the newly inserted row 37 with this week's data, matching the
# format/style of the surrounding rows (same as row 38's data before the shift).
$ws.Cells.Item(37, 1).Value = 1
$ws.Cells.Item(37, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(37, 4).Value = 45008
$ws.Cells.Item(37, 4).Style = $ws.Cells.Item(38, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
$ws.Cells.Item(37, 5).Value = 15
$ws.Cells.Item(37, 6).Value = 100112044
$ws.Cells.Item(37, 7).Value = "Perejil"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 250
$ws.Cells.Item(37, 11).Value = 2000
$ws.Cells.Item(37, 12).Value = 2500
$ws.Cells.Item(37, 13).Value = 2200
$ws.Cells.Item(37, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(37, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(37, 16).Value = 1100
$ws.Cells.Item(37, 17).Value = 2
$ws.Cells.Item(37, 18).Value = "Hortaliza"
